$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename milestone labels: "M" -> "M1" and "N" -> "M2"
$ws.Range("C5").Value = "M1"
$ws.Range("C6").Value = "M1 - A"
$ws.Range("C7").Value = "M1 - B"
$ws.Range("C8").Value = "M2"
$ws.Range("C9").Value = "M2 - C"
$ws.Range("C10").Value = "M2 - D"

# Add task fills for the milestone header rows, matching the fill already
# used for the sub-task cells beneath them (F6 and G7 for milestone 1,
# H9 and I10 for milestone 2).
$ws.Range("F6").Copy()
$ws.Range("F5:G5").PasteSpecial(-4122)

$ws.Range("H9").Copy()
$ws.Range("H8:I8").PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false
